$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# "Enero" row: plazo cell "Viernes 19 febrero 2021" -> "19 febrero 2021"
Replace-Text "Viernes 19 febrero 2021" "19 febrero 2021"

# "Febrero" row: plazo cell " Viernes 19 marzo 2021" -> " 19 marzo 2021"
Replace-Text " Viernes 19 marzo 2021" " 19 marzo 2021"

# "Marzo" row: plazo cell "Jueves 22" + " abril 2021" -> "22" + " abril 2021"
Replace-Text "Jueves 22 abril 2021" "22 abril 2021"

# "Abril" row: plazo cell "Lunes 24 de mayo 2021" -> "24 de mayo 2021"
Replace-Text "Lunes 24 de mayo 2021" "24 de mayo 2021"

# "Mayo" row: plazo cell "Lunes 21" + " de junio 2021" -> "21" + " de junio 2021"
Replace-Text "Lunes 21 de junio 2021" "21 de junio 2021"

# "Junio" row: plazo cell " Jueves" + " 22" + " julio 2021" -> " " + "22" + " julio 2021"
Replace-Text " Jueves 22 julio 2021" " 22 julio 2021"

# "Julio" row: plazo cell "Lunes 23" + " agosto 2021" -> "23" + " agosto 2021"
Replace-Text "Lunes 23 agosto 2021" "23 agosto 2021"

# "Agosto" row: plazo cell "Mi" + "é" + "rcoles 22" + " septiembre" + " 2021" -> "22" + " septiembre" + " 2021"
Replace-Text "Miércoles 22" "22"

# "Septiembre" row: plazo cell "Viernes 22 octubre 2021" -> "22 octubre 2021"
Replace-Text "Viernes 22 octubre 2021" "22 octubre 2021"

# "Octubre" row: plazo cell "Lunes 22 noviembre 2021" -> "22 noviembre 2021"
Replace-Text "Lunes 22 noviembre 2021" "22 noviembre 2021"

# "Noviembre" row: plazo cell "Miércoles 22 diciembre 2021" -> "22 diciembre 2021"
Replace-Text "Miércoles 22 diciembre 2021" "22 diciembre 2021"

# "Diciembre" row: plazo cell "Viernes 21" + " " + "enero " + "2022" -> "21" + " " + "enero " + "2022"
Replace-Text "Viernes 21 enero" "21 enero"
